# Apply "Generate Report for handback" update to localization-status.xlsx
#
# For each language sheet (zh-cn, de-de) the translated file has now been
# handed back in sync with en-US:
#   - Status / Handoff Reason move from "Ready for handoff" to
#     "Handed back: in sync with en-US" / "Include".
#   - The handed-back markdown file and its translated .xlf companion are
#     recorded as the "Latest Target File" (E2) and "Latest Handback File"
#     (F2) hyperlinks.
#   - "Latest Handback DateTime" (G2) is stamped with the handback time.
# The Overview sheet mirrors the same Status text, since it shares the
# underlying value with the language sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$mdFile    = "d655f31e-2332-4b20-a1ca-7bfaf702b797.md"
$mdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/e64a23f20aee8600bd755f230f5ffa202dad1ea7/e2e/d655f31e-2332-4b20-a1ca-7bfaf702b797.md"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/e64a23f20aee8600bd755f230f5ffa202dad1ea7/.localization-config"

# ---------------------------------------------------------------------
# zh-cn sheet (index 2)
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item(2)

$zhXlfFile        = "d655f31e-2332-4b20-a1ca-7bfaf702b797.6dc30c06dc380f80e049846d24089a724ab586d7.zh-cn.xlf"
$zhHandoffXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3e87a2d283a47d11481a85a69e6fadb52535d727/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/d655f31e-2332-4b20-a1ca-7bfaf702b797.6dc30c06dc380f80e049846d24089a724ab586d7.zh-cn.xlf"
$zhHandbackXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e64a23f20aee8600bd755f230f5ffa202dad1ea7/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/d655f31e-2332-4b20-a1ca-7bfaf702b797.6dc30c06dc380f80e049846d24089a724ab586d7.zh-cn.xlf"

$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("G2").Value = "2016-01-17 14:31:08"
$wsZh.Range("H2").Value = "Include"

# Rebuild the hyperlinks collection so the new "Latest Target File" (E2) and
# "Latest Handback File" (F2) links take their place, in sheet order,
# alongside the pre-existing ones.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhHandoffXlfUrl, [Type]::Missing, [Type]::Missing, $zhXlfFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhHandbackXlfUrl, [Type]::Missing, [Type]::Missing, $zhXlfFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $configUrl, [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet (index 3)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item(3)

$deXlfFile        = "d655f31e-2332-4b20-a1ca-7bfaf702b797.6dc30c06dc380f80e049846d24089a724ab586d7.de-de.xlf"
$deHandoffXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/97114c3bb492eb88cbb1b0cb771b6f8eab71c8ac/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/d655f31e-2332-4b20-a1ca-7bfaf702b797.6dc30c06dc380f80e049846d24089a724ab586d7.de-de.xlf"
$deHandbackXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e64a23f20aee8600bd755f230f5ffa202dad1ea7/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/d655f31e-2332-4b20-a1ca-7bfaf702b797.6dc30c06dc380f80e049846d24089a724ab586d7.de-de.xlf"

$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("G2").Value = "2016-01-17 14:31:26"
$wsDe.Range("H2").Value = "Include"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deHandoffXlfUrl, [Type]::Missing, [Type]::Missing, $deXlfFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deHandbackXlfUrl, [Type]::Missing, [Type]::Missing, $deXlfFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $configUrl, [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Overview sheet status column mirrors the language sheets' Status cell.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

$wb.Save()
